$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-15: new "Actualizado" timestamp
$ws.Range("D2:D15").Value = 44248.53816747577

# Rows 16-29: shift down from what rows 2-15 used to hold
$ws.Range("D16:D29").Value = 44248.51686203704

# Rows 30-43: shift down from what rows 16-29 used to hold
$ws.Range("D30:D43").Value = 44248.49556476852
